$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.Value = "'" + $text
    $rng.Style = $origStyle
}

Set-TextValue "C11" "18.19"

Set-TextValue "B30" "71.48"
Set-TextValue "C30" "18.48"
Set-TextValue "D30" "89.96"

Set-TextValue "B31" "16.49"
Set-TextValue "C31" "41.02"
Set-TextValue "D31" "57.51"

Set-TextValue "B33" "79.41"
Set-TextValue "C33" "20.53"
Set-TextValue "D33" "99.94"
